# Updated cryptos list on Wed Apr 24 17:55:39 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.875.26"
$ws.Range("E2").Value = "  -2.81%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.180.81"
$ws.Range("E3").Value = "  -1.73%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'602.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.25%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'152.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.45%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.178.16"
$ws.Range("E8").Value = "  -1.71%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -3.32%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -4.81%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -1.89%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.479"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.59%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000262"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.25%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'37.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.89%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.682.93"
$ws.Range("E15").Value = "  -2.31%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "64.943.98"
$ws.Range("E16").Value = "  -2.78%  "

# Row 17 - WrappedEther -> TRON
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.114"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.61%  "

# Row 18 - TRON -> WrappedEther
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.156.80"
$ws.Range("E18").Value = "  -2.49%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -4.08%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'485.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.34%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'14.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.95%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.55%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "'7.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.47%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'14.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.36%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'85.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "

# Row 26 - Dai
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "

# Row 27 - PancakeSwap
$ws.Range("D27").Value = "'2.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.37%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'8.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.99%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  -4.42%  "

# Row 30 - NEARProtocol
$ws.Range("D30").Value = "'7.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.27%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +0.32%  "

# Row 32 - Stacks
$ws.Range("D32").Value = "'2.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.86%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'27.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.40%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  -0.04%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -5.59%  "

# Row 36 - Filecoin
$ws.Range("D36").Value = "'6.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.51%  "

# Row 37 - dogwifhat -> OKB
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'54.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.54%  "

# Row 38 - OKB -> dogwifhat
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.25%  "

# Row 39 - PEPE
$ws.Range("D39").Value = "0.0₃0750"
$ws.Range("E39").Value = "  -3.00%  "

# Row 40 - Bittensor
$ws.Range("D40").Value = "'460.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.57%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  -2.00%  "

# Row 42 - VeChain
$ws.Range("D42").Value = "'0.0407"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.80%  "

# Row 43 - Cosmos
$ws.Range("D43").Value = "'8.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.33%  "

# Row 44 - Fetch.AI
$ws.Range("E44").Value = "  +0.42%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.922.22"
$ws.Range("E45").Value = "  +0.14%  "

# Row 46 - TheGraph
$ws.Range("D46").Value = "'0.278"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.37%  "

# Row 47 - InjectiveProtocol
$ws.Range("D47").Value = "'27.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.54%  "

# Row 48 - USDe
$ws.Range("E48").Value = "  +0.02%  "

# Row 49 - ThetaToken
$ws.Range("D49").Value = "'2.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.13%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  -0.10%  "

# Row 51 - Monero
$ws.Range("D51").Value = "'120.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.12%  "

